$wb = $excel.ActiveWorkbook

# --- Sheet "Alluvial for Mapping": fix misspelled watershed name ---
$wsMap = $wb.Worksheets.Item("Alluvial for Mapping")
for ($r = 8; $r -le 57; $r++) {
    $cell = $wsMap.Cells.Item($r, 16)  # column P = Watershed
    if ($cell.Text -eq "Mortendad") {
        $cell.Value = "Mortandad"
    }
}

# --- Sheet "Alluvial Exhibit": fixes + restructuring ---
$wsExh = $wb.Worksheets.Item("Alluvial Exhibit")

# Widen column C slightly
$wsExh.Columns.Item(3).ColumnWidth = 12.14

# Fix misspelled canyon heading
$wsExh.Range("A10").Value = "Mortandad Canyon"

# Split "Los Alamos and Pajarito Canyons" heading into two separate sections
$wsExh.Range("A61").Value = "Los Alamos Canyon"

# Insert a new heading row for Pajarito Canyon before the former row 74 (18-BG-4)
$wsExh.Rows.Item(74).Insert()
$wsExh.Range("A74:H74").Merge()
$wsExh.Range("A61:H61").Copy()
$wsExh.Range("A74").PasteSpecial(-4122)
$wsExh.Range("A74").Value = "Pajarito Canyon"
